# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows at the top of the Kiwi data block (rows 364-367),
# pushing the existing data (old rows 364-429) down to rows 368-433.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 364 (formatting copied from the surrounding rows)
$ws.Range("A364:A367").EntireRow.Insert()

# Values shared by all four new rows
$mercadoId   = 9
$mercado     = 'Vega Central Mapocho de Santiago'
$region      = 'Metropolitana'
$fecha       = 44504
$codreg      = 13
$tipo        = 'Fruta'
$productoId  = 100101
$producto    = 'Berries'
$categoriaId = 100101007
$categoria   = 'Kiwi'
$variedad    = 'Hayward'
$unidad      = '$/caja 18 kilos'
$kgUnidad    = 18

# Row 364: Especial
$r = 364
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = 'Especial'
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 19800
$ws.Cells.Item($r, 15).Value = 19800
$ws.Cells.Item($r, 16).Value = 19800
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item($r, 19).Value = 1100
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 365: Extra (doble especial)
$r = 365
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = 'Extra (doble especial)'
$ws.Cells.Item($r, 13).Value = 150
$ws.Cells.Item($r, 14).Value = 21600
$ws.Cells.Item($r, 15).Value = 21600
$ws.Cells.Item($r, 16).Value = 21600
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item($r, 19).Value = 1200
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 366: Primera
$r = 366
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = 'Primera'
$ws.Cells.Item($r, 13).Value = 220
$ws.Cells.Item($r, 14).Value = 18000
$ws.Cells.Item($r, 15).Value = 18000
$ws.Cells.Item($r, 16).Value = 18000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item($r, 19).Value = 1000
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 367: Segunda
$r = 367
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = 'Segunda'
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 6500
$ws.Cells.Item($r, 15).Value = 6500
$ws.Cells.Item($r, 16).Value = 6500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item($r, 19).Value = 361
$ws.Cells.Item($r, 20).Value = $kgUnidad

Write-Host "Applied 4 new rows; dimension should now be A1:T433"
